$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows above the existing row 2 (old row2 becomes row5) ---
$ws.Rows("2:4").Insert()

# --- Fill the 3 newly inserted rows (2,3,4) ---
$ws.Range("A2").Value = 43444
$ws.Range("B2").Value = "Absent"
$ws.Range("C2").Value = "Recherche de QCM et de templates"

$ws.Range("A3").Value = 43446
$ws.Range("B3").Value = "Analyse des sites et compléter la grille d'évaluation"
$ws.Range("C3").Value = "Analyse des sites et compléter la grille d'évaluation"

$ws.Range("A4").Value = 43451
$ws.Range("B4").Value = "Conception du Moodboard, wireframe balsamiq"
$ws.Range("C4").Value = "Synthèse de l'analyse, analyse des 2 derniers templates du sites de façon moins détaillée "

# --- Update the text of the (shifted) original data row, now row 5 ---
$ws.Range("C5").Value = "Creation maquette, création des templates sur photoshop, création logo"

# --- Append two new rows (6 and 7) after the existing data ---
$ws.Range("A6").Value = 43109
$ws.Range("C6").Value = "Finition maquette et repartition du travail"

$ws.Range("A7").Value = 43114
$ws.Range("C7").Value = "Création d'un template complet"

# --- Styling: center-align header/column A, and re-apply centered date format ---
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").NumberFormat = "mm-dd-yy"

$ws.Range("A2").Copy()
$ws.Range("A3:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the selection shown on the sheet (matches final diff state) ---
$ws.Range("C9").Select()
